# Applies the "Update data: 2025-11-20 09:18" edit described by the diff.
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Metadata sheet: bump the "Last Updated" timestamp (A2)
# ---------------------------------------------------------------------
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("A2").Value = "20 Nov 2025, 09:18 AM"

# ---------------------------------------------------------------------
# 2. Top Gainers sheet: update the Weekly value (D75) for STUDDS
# ---------------------------------------------------------------------
$wsGainers = $wb.Worksheets.Item("Top Gainers")
$wsGainers.Range("D75").Value = 5.7036

# ---------------------------------------------------------------------
# 3. Industry Analysis sheet: the "personal care - indian" row's Weekly
#    value (column D) dropped from 3.4188 to 1.9824, which moves that
#    row from the top of this descending-sorted block (row 11) down to
#    the bottom (row 17); rows 12-17 each shift up by one row. Also
#    row 65's Latest value (C65) is updated independently.
# ---------------------------------------------------------------------
$wsIndustry = $wb.Worksheets.Item("Industry Analysis")

# New contents for rows 11-17 (columns B..K), after the re-sort
$rows = @(
    @{ Row = 11; B = "automobiles - motorcycles / mopeds";            C = 0.5837;  D = 3.2925; E = -0.305;   F = 47.6884;  G = 16.9168;  H = 86.6409;  I = 29.2784;  J = 41.8958;  K = 31.174 }
    @{ Row = 12; B = "cement products";                                C = -1.9048; D = 2.7124; E = 0.6843;   F = -51.4609; G = 29.4692;  H = 13.8194;  I = 150.4245; J = 114.8705; K = 79.1866 }
    @{ Row = 13; B = "dry cells";                                      C = 0;       D = 2.6598; E = -10.5949; F = -13.3656; G = 13.2601;  H = -4.4972;  I = 38.329;   J = 71.5429;  K = 37.2103 }
    @{ Row = 14; B = "bearings";                                       C = 0.6332;  D = 2.4033; E = 7.3777;   F = 5.1126;   G = -16.7188; H = 120.0131; I = 58.3623;  J = 32.445;   K = 17.366 }
    @{ Row = 15; B = "pesticides / agrochemicals - multinational";     C = 0.2655;  D = 2.0955; E = -6.4407;  F = -17.6718; G = -0.0162;  H = 12.4448;  I = -8.7623;  J = 11.1706;  K = 7.1425 }
    @{ Row = 16; B = "auto ancillaries";                                C = 0.4621;  D = 1.9988; E = 4.8424;   F = 6.5208;   G = 46.2142;  H = 36.882;   I = 67.149;   J = 28.2354;  K = 11.7715 }
    @{ Row = 17; B = "personal care - indian";                         C = 0.1866;  D = 1.9824; E = -0.3932;  F = 39.84;    G = -10.0911; H = 32.0505;  I = -7.6062;  J = -2.178;   K = -3.3673 }
)

foreach ($r in $rows) {
    $rowNum = $r.Row
    $wsIndustry.Range("B$rowNum").Value = $r.B
    $wsIndustry.Range("C$rowNum").Value = $r.C
    $wsIndustry.Range("D$rowNum").Value = $r.D
    $wsIndustry.Range("E$rowNum").Value = $r.E
    $wsIndustry.Range("F$rowNum").Value = $r.F
    $wsIndustry.Range("G$rowNum").Value = $r.G
    $wsIndustry.Range("H$rowNum").Value = $r.H
    $wsIndustry.Range("I$rowNum").Value = $r.I
    $wsIndustry.Range("J$rowNum").Value = $r.J
    $wsIndustry.Range("K$rowNum").Value = $r.K
}

# Independent value update on row 65 (automobiles - lcvs / hcvs, Latest column)
$wsIndustry.Range("C65").Value = -1.0144
